# Deploying to gh-pages: add 2021 data column (T) to the "number of children
# not attending school" table, mirroring the existing 2020 column (S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlRight = -4152

# ---------------------------------------------------------------------------
# Structural rows: copy S-column formatting verbatim (no alignment change).
# ---------------------------------------------------------------------------
$structural = @(
    @{Row=3;  Val=$null},
    @{Row=4;  Val=2021},
    @{Row=5;  Val=$null},
    @{Row=23; Val=$null}
)

foreach ($item in $structural) {
    $src = "S" + $item.Row
    $dst = "T" + $item.Row
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial($xlPasteFormats)
    if ($null -ne $item.Val) {
        $ws.Range($dst).Value = $item.Val
    }
}

# ---------------------------------------------------------------------------
# Data rows: copy S-column formatting, then force right-aligned text (the
# new 2021 column is right-aligned throughout the data area), then write the
# 2021 value. "-" values reuse the existing shared string automatically.
# ---------------------------------------------------------------------------
$data = @(
    @{Row=6;  Val=1466},
    @{Row=7;  Val=$null},
    @{Row=8;  Val=76},
    @{Row=9;  Val=15},
    @{Row=10; Val=1},
    @{Row=11; Val=188},
    @{Row=12; Val=22},
    @{Row=13; Val=15},
    @{Row=14; Val="-"},
    @{Row=15; Val="-"},
    @{Row=16; Val=112},
    @{Row=17; Val="-"},
    @{Row=18; Val=6},
    @{Row=19; Val="-"},
    @{Row=20; Val=29},
    @{Row=21; Val=1002},
    @{Row=22; Val="-"},
    @{Row=24; Val=1029},
    @{Row=25; Val=$null},
    @{Row=26; Val=51},
    @{Row=27; Val=4},
    @{Row=28; Val="-"},
    @{Row=29; Val=127},
    @{Row=30; Val=14},
    @{Row=31; Val=12},
    @{Row=32; Val="-"},
    @{Row=33; Val="-"},
    @{Row=34; Val=70},
    @{Row=35; Val="-"},
    @{Row=36; Val=3},
    @{Row=37; Val="-"},
    @{Row=38; Val=16},
    @{Row=39; Val=732},
    @{Row=40; Val="-"}
)

foreach ($item in $data) {
    $src = "S" + $item.Row
    $dst = "T" + $item.Row
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial($xlPasteFormats)
    $ws.Range($dst).HorizontalAlignment = $xlRight
    if ($null -ne $item.Val) {
        $ws.Range($dst).Value = $item.Val
    }
}

# ---------------------------------------------------------------------------
# Final selection matches the author's saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("T3").Select()
